# OPSK-1387 NeLS sample extraction should run in "overwrite" mode
#
# Appends a 6th sample row ("Sample-5" / "Lib-1" / "e") to the "Samples"
# sheet, reusing the existing "Lib-1" shared string for column B (library
# id) just like the real extraction run did, and moves the active
# selection to the newly written C6 cell (mirrors the author's workbook
# state after entering the row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Sample-5"
$ws.Range("B6").Value = "Lib-1"
$ws.Range("C6").Value = "e"

$ws.Range("C6").Select() | Out-Null
